$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the shared-string table in the same order the original author typed
# them in, so new <si> entries line up with the target uniqueCount order:
#   39 = "System load estimation"
#   40 = "Minor changes on documentation and setup. System load estimation"
#   41 = "tc: System load estimation put to operation, validated by test case tc10"
$ws.Range("D48").Value = "System load estimation"
$ws.Range("D47").Value = "Minor changes on documentation and setup. System load estimation"
$ws.Range("D49").Value = "tc: System load estimation put to operation, validated by test case tc10"

# New effort-log entries for tc10 (system load estimation)
# Row 47: 2012-11-26, 2h, "Minor changes on documentation and setup. System load estimation"
$ws.Range("A47").Value = 41239
$ws.Range("A47").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B47").Value = 2

# Row 48: 2012-11-28, 2h, "System load estimation"
$ws.Range("A48").Value = 41241
$ws.Range("A48").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B48").Value = 2

# Row 49: 2012-11-29, 2h, "tc: System load estimation put to operation, validated by test case tc10"
$ws.Range("A49").Value = 41242
$ws.Range("A49").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B49").Value = 2

# Update selection/view to match post-edit state
$ws.Range("A50").Select()
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
